$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix workbook-level absolute path metadata (folder name swap) ---
# (not directly settable via the object model; handled separately if supported)

# --- data correction: Khalil's math score 56 -> 57 (row 4 in original layout) ---
$ws.Range("D4").Value = 57

# --- header for new Rank column ---
$ws.Range("H1").Value = "Rank"

# --- replace average / grade formulas and add rank formulas, row by row ---
$ws.Range("F2").Formula = "=AVERAGE(B2:E2)"
$ws.Range("G2").Formula = "=IF(MIN(B2:E2)<40,""F"",IF(F2>=80,""A+"",IF(F2>=70,""A-"",IF(F2>=60,""B"",""C""))))"
$ws.Range("H2").Formula = "=RANK(F2,F`$2:F`$5,0)"

$ws.Range("F3").Formula = "=AVERAGE(B3:E3)"
$ws.Range("G3").Formula = "=IF(MIN(B3:E3)<40,""F"",IF(F3>=80,""A+"",IF(F3>=70,""A-"",IF(F3>=60,""B"",""C""))))"
$ws.Range("H3").Formula = "=RANK(F3,F`$2:F`$5,0)"

$ws.Range("F4").Formula = "=AVERAGE(B4:E4)"
$ws.Range("G4").Formula = "=IF(MIN(B4:E4)<40,""F"",IF(F4>=80,""A+"",IF(F4>=70,""A-"",IF(F4>=60,""B"",""C""))))"
$ws.Range("H4").Formula = "=RANK(F4,F`$2:F`$5,0)"

$ws.Range("F5").Formula = "=AVERAGE(B5:E5)"
$ws.Range("G5").Formula = "=IF(MIN(B5:E5)<40,""F"",IF(F5>=80,""A+"",IF(F5>=70,""A-"",IF(F5>=60,""B"",""C""))))"
$ws.Range("H5").Formula = "=RANK(F5,F`$2:F`$5,0)"

# --- sort the student table by the new Rank column (ascending) ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("H2:H5"))
$ws.Sort.SetRange($ws.Range("A1:H5"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# --- apply then remove AutoFilter, leaving the hidden _FilterDatabase name behind ---
$ws.Range("A1:H5").AutoFilter()
$ws.AutoFilterMode = $false
$n = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$H`$5")
$n.Visible = $false

# --- selection left on H4 after the edits ---
$ws.Range("H4").Select()

Write-Host "done"
